$d = $word.ActiveDocument

# 1. Increase the "Talles" title font size from 20pt (sz=40) to 28pt (sz=56).
$p1 = $d.Paragraphs.Item(1)
$p1.Range.Font.Size = 28

# 2. Move the "_GoBack" bookmark from the "Como encontrar tu talle!" paragraph
#    to the end of the "Talles" paragraph (right after the run, before the
#    paragraph mark).
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

$r = $p1.Range.Duplicate
$null = $r.MoveEnd(1, -1)
$null = $r.Collapse(0)
$null = $r.InsertAfter("X")
$insertPos = $r.Start
$rb = $d.Range($insertPos, $insertPos)
$null = $d.Bookmarks.Add("_GoBack", $rb)
$rx = $d.Range($insertPos, $insertPos + 1)
$null = $rx.Delete()
